$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O (15, "بانك"/Bank) by shifting every
# populated header cell from AK (37) down to O (15) one column to the
# right (iterate right-to-left so sources aren't clobbered before they are
# read). Range.Copy brings the cell style (s attribute) along with it, so
# this reproduces an "insert column" without touching the <cols> ranges.
for ($c = 37; $c -ge 15; $c--) {
    $src = $ws.Cells.Item(1, $c)
    $dst = $ws.Cells.Item(1, $c + 1)
    $src.Copy($dst)
}

# Populate the freed-up column O with the new header text.
$ws.Cells.Item(1, 15).Value2 = "شماره شبا"

# Reflect the new used range in the filter-database defined name.
$wb.Names("_xlnm._FilterDatabase").RefersTo = '=Sheet1!$A$1:$WXA$1'

# Update the view state to match where the user left the selection.
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("O1").Select()
